$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.8419516666666667
$ws.Range("H2").Value = 2.525855
$ws.Range("I2").Value = 0.005772335854373203
$ws.Range("J2").Value = 0.005772335854373203
$ws.Range("M2").Value = 14.129345
$ws.Range("N2").Value = 42.388035
$ws.Range("O2").Value = 0.3414817166893976
$ws.Range("P2").Value = 0.3414817166893976
$ws.Range("Q2").Value = 11.89622557165833
$ws.Range("R2").Value = 107.066030144925
$ws.Range("S2").Value = 0.001971147156859122
$ws.Range("T2").Value = 0.001971147156859122

$ws.Range("G3").Value = 0.8419516666666667
$ws.Range("H3").Value = 2.525855
$ws.Range("I3").Value = 0.005772335854373203
$ws.Range("J3").Value = 0.005772335854373203
$ws.Range("O3").Value = 0.3553528814026711
$ws.Range("P3").Value = 0.3553528814026711
$ws.Range("Q3").Value = 12.37945643382722
$ws.Range("R3").Value = 111.415107904445
$ws.Range("S3").Value = 0.002051216178275467
$ws.Range("T3").Value = 0.002051216178275467

$ws.Range("G4").Value = 0.8419516666666667
$ws.Range("H4").Value = 2.525855
$ws.Range("I4").Value = 0.005772335854373203
$ws.Range("J4").Value = 0.005772335854373203
$ws.Range("O4").Value = 0.3031654019079313
$ws.Range("P4").Value = 0.3031654019079312
$ws.Range("Q4").Value = 10.56139708322833
$ws.Range("R4").Value = 95.052573749055
$ws.Range("S4").Value = 0.001749972519238614
$ws.Range("T4").Value = 0.001749972519238614

$ws.Range("I5").Value = 0.8411037170617888
$ws.Range("J5").Value = 0.8411037170617888
$ws.Range("M5").Value = 14.129345
$ws.Range("N5").Value = 42.388035
$ws.Range("O5").Value = 0.3414817166893976
$ws.Range("P5").Value = 0.3414817166893976
$ws.Range("Q5").Value = 1733.433362119197
$ws.Range("R5").Value = 15600.90025907277
$ws.Range("S5").Value = 0.287221541216093
$ws.Range("T5").Value = 0.287221541216093

$ws.Range("I6").Value = 0.8411037170617888
$ws.Range("J6").Value = 0.8411037170617888
$ws.Range("O6").Value = 0.3553528814026711
$ws.Range("P6").Value = 0.3553528814026711
$ws.Range("S6").Value = 0.2988886294164037
$ws.Range("T6").Value = 0.2988886294164037

$ws.Range("I7").Value = 0.8411037170617888
$ws.Range("J7").Value = 0.8411037170617888
$ws.Range("O7").Value = 0.3031654019079313
$ws.Range("P7").Value = 0.3031654019079312
$ws.Range("S7").Value = 0.2549935464292921
$ws.Range("T7").Value = 0.2549935464292921

$ws.Range("I8").Value = 0.1531239470838381
$ws.Range("J8").Value = 0.1531239470838381
$ws.Range("M8").Value = 14.129345
$ws.Range("N8").Value = 42.388035
$ws.Range("O8").Value = 0.3414817166893976
$ws.Range("P8").Value = 0.3414817166893976
$ws.Range("Q8").Value = 315.57363619997
$ws.Range("R8").Value = 2840.16272579973
$ws.Range("S8").Value = 0.0522890283164455
$ws.Range("T8").Value = 0.0522890283164455

$ws.Range("I9").Value = 0.1531239470838381
$ws.Range("J9").Value = 0.1531239470838381
$ws.Range("O9").Value = 0.3553528814026711
$ws.Range("P9").Value = 0.3553528814026711
$ws.Range("S9").Value = 0.05441303580799201
$ws.Range("T9").Value = 0.05441303580799201

$ws.Range("I10").Value = 0.1531239470838381
$ws.Range("J10").Value = 0.1531239470838381
$ws.Range("O10").Value = 0.3031654019079313
$ws.Range("P10").Value = 0.3031654019079312
$ws.Range("S10").Value = 0.04642188295940058
$ws.Range("T10").Value = 0.04642188295940057
